# Update cryptocurrency price/volume data per Thu Jul 18 08:34:26 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.756.96"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.440.83"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "572.05"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "158.83"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.440.17"
$ws.Range("E8").Value = "  -1.67%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.572"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -6.25%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.121"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.81%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.441"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "4.034.61"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E14").Value = "  -0.59%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.64"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.88%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000175"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -9.79%  "
$ws.Range("D17").Value = "64.785.07"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "3.468.35"
$ws.Range("E18").Value = "  -0.25%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.23"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.78%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.82"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.94%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "377.61"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("E24").Value = "  -0.01%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "71.91"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  -0.66%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.84"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("E28").Value = "  -1.07%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.47"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.38%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.07"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.78%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.67%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.12"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.47%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.97"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -1.93%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "161.05"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.49%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.88"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").Value = "2.894.43"
$ws.Range("E38").Value = "  -4.00%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0747"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("E40").Value = "  +2.77%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.16"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.76%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.51"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  +0.46%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.780"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "25.98"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0310"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "318.18"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.08"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.82%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.47"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.63%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.842"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.41%  "
